$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - plain numeric header values 1..6 in F1:K1
$ws.Range("F1").Value = 1
$ws.Range("G1").Value = 2
$ws.Range("H1").Value = 3
$ws.Range("I1").Value = 4
$ws.Range("J1").Value = 5
$ws.Range("K1").Value = 6

# Row 2 - COUNTIF(..., "Diabetes") across columns B..G, plus row total in M2
$ws.Range("F2").Formula = '=COUNTIF(B:B, "Diabetes")'
$ws.Range("G2:K2").Formula = '=COUNTIF(C:C, "Diabetes")'
$ws.Range("M2").Formula = '=SUM(F2:K2)'

# Row 3 - COUNTIF(..., "*") across columns B..G, plus row total in M3
$ws.Range("F3").Formula = '=COUNTIF(B:B, "*")'
$ws.Range("G3:K3").Formula = '=COUNTIF(C:C, "*")'
$ws.Range("M3").Formula = '=SUM(F3:K3)'

# Row 4 - ratio of row2/row3, formatted as percentage; L4 stays blank (styled only)
$ws.Range("F4").Formula = '=F2/F3'
$ws.Range("G4:M4").Formula = '=G2/G3'
$ws.Range("F4:M4").NumberFormat = "0.00%"
$ws.Range("L4").ClearContents()

# Restore the view state (selection/zoom) as left by the author
$null = $ws.Range("I10").Select()
$ws.Application.ActiveWindow.Zoom = 172
